$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 is new; populate Year/Month/Day for it
$ws.Cells.Item(16, 1).Value = 2026
$ws.Cells.Item(16, 2).Value = 2
$ws.Cells.Item(16, 3).Value = 19

# Hourly load values for rows 11-16, columns D (4) through AA (27)
$data = New-Object 'object[,]' 6,24
$data[0,0] = 4678
$data[0,1] = 4466
$data[0,2] = 4380
$data[0,3] = 4350
$data[0,4] = 4343
$data[0,5] = 4374
$data[0,6] = 4522
$data[0,7] = 4798
$data[0,8] = 5288
$data[0,9] = 5586
$data[0,10] = 5661
$data[0,11] = 5574
$data[0,12] = 5463
$data[0,13] = 5368
$data[0,14] = 5223
$data[0,15] = 5152
$data[0,16] = 5250
$data[0,17] = 5533
$data[0,18] = 5816
$data[0,19] = 5747
$data[0,20] = 5472
$data[0,21] = 5123
$data[0,22] = 4980
$data[0,23] = 4716
$data[1,0] = 4368
$data[1,1] = 4225
$data[1,2] = 4160
$data[1,3] = 4099
$data[1,4] = 4079
$data[1,5] = 4125
$data[1,6] = 4285
$data[1,7] = 4630
$data[1,8] = 5070
$data[1,9] = 5394
$data[1,10] = 5451
$data[1,11] = 5392
$data[1,12] = 5309
$data[1,13] = 5252
$data[1,14] = 5109
$data[1,15] = 5061
$data[1,16] = 5229
$data[1,17] = 5513
$data[1,18] = 5865
$data[1,19] = 5867
$data[1,20] = 5677
$data[1,21] = 5389
$data[1,22] = 5161
$data[1,23] = 4836
$data[2,0] = 4501
$data[2,1] = 4303
$data[2,2] = 4255
$data[2,3] = 4203
$data[2,4] = 4246
$data[2,5] = 4427
$data[2,6] = 4935
$data[2,7] = 5476
$data[2,8] = 5919
$data[2,9] = 6130
$data[2,10] = 5894
$data[2,11] = 5900
$data[2,12] = 5558
$data[2,13] = 5555
$data[2,14] = 5538
$data[2,15] = 5602
$data[2,16] = 5677
$data[2,17] = 5866
$data[2,18] = 6267
$data[2,19] = 6285
$data[2,20] = 6099
$data[2,21] = 5835
$data[2,22] = 5602
$data[2,23] = 5197
$data[3,0] = 4766
$data[3,1] = 4572
$data[3,2] = 4491
$data[3,3] = 4489
$data[3,4] = 4494
$data[3,5] = 4662
$data[3,6] = 5188
$data[3,7] = 5755
$data[3,8] = 6203
$data[3,9] = 6343
$data[3,10] = 6233
$data[3,11] = 6099
$data[3,12] = 5957
$data[3,13] = 5837
$data[3,14] = 5696
$data[3,15] = 5683
$data[3,16] = 5708
$data[3,17] = 5957
$data[3,18] = 6492
$data[3,19] = 6488
$data[3,20] = 6285
$data[3,21] = 6015
$data[3,22] = 5781
$data[3,23] = 5415
$data[4,0] = 4958
$data[4,1] = 4763
$data[4,2] = 4681
$data[4,3] = 4679
$data[4,4] = 4684
$data[4,5] = 4853
$data[4,6] = 5380
$data[4,7] = 5948
$data[4,8] = 6397
$data[4,9] = 6538
$data[4,10] = 6417
$data[4,11] = 6271
$data[4,12] = 6115
$data[4,13] = 5985
$data[4,14] = 5830
$data[4,15] = 5816
$data[4,16] = 5843
$data[4,17] = 6109
$data[4,18] = 6680
$data[4,19] = 6676
$data[4,20] = 6476
$data[4,21] = 6209
$data[4,22] = 5977
$data[4,23] = 5615
$data[5,0] = 4950
$data[5,1] = 4755
$data[5,2] = 4673
$data[5,3] = 4672
$data[5,4] = 4677
$data[5,5] = 4845
$data[5,6] = 5372
$data[5,7] = 5940
$data[5,8] = 6389
$data[5,9] = 6531
$data[5,10] = 6410
$data[5,11] = 6264
$data[5,12] = 6109
$data[5,13] = 5979
$data[5,14] = 5825
$data[5,15] = 5811
$data[5,16] = 5838
$data[5,17] = 6103
$data[5,18] = 6673
$data[5,19] = 6669
$data[5,20] = 6468
$data[5,21] = 6201
$data[5,22] = 5969
$data[5,23] = 5607

$ws.Range("D11:AA16").Value = $data
